# Generate Report for Handoff
# Adds a new "ed1873f5-da6e-4f21-8056-7c55f07ef4ec" handoff row to each of the
# three worksheets (Overview, zh-cn, de-de), mirroring the existing
# "27e84707-347a-4dd8-97db-1852d3153c02" row that is already present on row 2.

$wb = $excel.ActiveWorkbook

$sourceMd       = "ed1873f5-da6e-4f21-8056-7c55f07ef4ec.md"
$zhXlf          = "ed1873f5-da6e-4f21-8056-7c55f07ef4ec.3f243414d1a18aa75cdbdbe1caf15cfd73b76cff.zh-cn.xlf"
$deXlf          = "ed1873f5-da6e-4f21-8056-7c55f07ef4ec.3f243414d1a18aa75cdbdbe1caf15cfd73b76cff.de-de.xlf"

$mdCommit       = "c2c7f8a9d7b3b6b1a9e2f5c4d6e8b1a3c5d7e9f1"
$zhCommit       = "3f243414d1a18aa75cdbdbe1caf15cfd73b76cff"
$deCommit       = "3f243414d1a18aa75cdbdbe1caf15cfd73b76cff"

$mdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/$sourceMd"
$zhUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf"
$deUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf"

$handoffDatetime = "2016-03-21 12:31:33"
$zhHandoffFileDt = "2016-03-21 12:31:27"
$deHandoffFileDt = "2016-03-21 12:31:33"
$epoch           = "0001-01-01 00:00:00"
$readyStatus     = "Ready for handoff"
$dateFmt         = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "Overview" -> row 3
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $mdUrl, [Type]::Missing, [Type]::Missing, $sourceMd)
$wsOverview.Range("B3").Value = $readyStatus
$wsOverview.Range("C3").Value = $readyStatus
$wsOverview.Range("D3").Value = $handoffDatetime
$wsOverview.Range("D3").NumberFormat = $dateFmt

# ---------------------------------------------------------------------------
# Sheet "zh-cn" -> row 3
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $mdUrl, [Type]::Missing, [Type]::Missing, $sourceMd)
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = $readyStatus
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $zhUrl, [Type]::Missing, [Type]::Missing, $zhXlf)
$wsZh.Range("E3").Value = $zhHandoffFileDt
$wsZh.Range("E3").NumberFormat = $dateFmt
$wsZh.Range("H3").Value = $epoch
$wsZh.Range("H3").NumberFormat = $dateFmt
$wsZh.Range("J3").Value = "Include"

# ---------------------------------------------------------------------------
# Sheet "de-de" -> row 3
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $mdUrl, [Type]::Missing, [Type]::Missing, $sourceMd)
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = $readyStatus
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $deUrl, [Type]::Missing, [Type]::Missing, $deXlf)
$wsDe.Range("E3").Value = $deHandoffFileDt
$wsDe.Range("E3").NumberFormat = $dateFmt
$wsDe.Range("H3").Value = $epoch
$wsDe.Range("H3").NumberFormat = $dateFmt
$wsDe.Range("J3").Value = "Include"

Write-Host "Handoff report row added for $sourceMd"
